# Fix unit bonuses in the "Units" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# H column = life-bonus-add, I column = life-bonus-mult, U column = dps-bonus-add
$ws.Range("U4").Value = 10        # Encoder: dps-bonus-add 8 -> 10
$ws.Range("H9").Value = 5280      # Archive: life-bonus-add 424 -> 5280
$ws.Range("I12").Value = 1.2      # Matriarch: life-bonus-mult 2.05 -> 1.2
$ws.Range("H21").Value = 1300     # SW3-NTZ WELDTECH RELAY: life-bonus-add 0 -> 1300
$ws.Range("H33").Value = 4500     # TOR-N0 UPHOLDER: life-bonus-add 4000 -> 4500
$ws.Range("H44").Value = 5263     # Mecha Host: life-bonus-add 5000 -> 5263
$ws.Range("H45").Value = 3158     # Primal Host: life-bonus-add 3000 -> 3158
$ws.Range("U56").Value = 60       # Archangel: dps-bonus-add 42.86 -> 60
$ws.Range("H63").Value = 5325     # Theos: life-bonus-add 3195 -> 5325
$ws.Range("U63").Value = 574      # Theos: dps-bonus-add 344 -> 574
$ws.Range("H83").Value = 5764     # Soul of Legends: life-bonus-add 0 -> 5764
$ws.Range("H116").Value = 1520    # YGGDRASIL: life-bonus-add 1560 -> 1520
$ws.Range("H117").Value = 5700    # T.R.E.E. of Life: life-bonus-add 4881 -> 5700
